$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 69556
$ws.Range("E9").Value = 191321361

$ws.Range("C10").Value = 278093
$ws.Range("E10").Value = 1751026417

$ws.Range("C14").Value = 119174
$ws.Range("D14").Value = 25290
$ws.Range("E14").Value = 379158686

$ws.Range("C19").Value = 108879
$ws.Range("E19").Value = 344360661

$ws.Range("C46").Value = 55743
$ws.Range("E46").Value = 174066006

$ws.Range("C64").Value = 10534
$ws.Range("E64").Value = 25235618

$ws.Range("C81").Value = 26152
$ws.Range("E81").Value = 164983877

$ws.Range("C99").Value = 136529
$ws.Range("E99").Value = 862699431

$ws.Range("C102").Value = 13703
$ws.Range("E102").Value = 28894138

$ws.Range("C104").Value = 22086
$ws.Range("E104").Value = 84715947

$ws.Range("C108").Value = 45968
$ws.Range("E108").Value = 145884559

$ws.Range("C111").Value = 5908
$ws.Range("E111").Value = 11622945

$ws.Range("C115").Value = 17127
$ws.Range("E115").Value = 37513971

$ws.Range("C150").Value = 94995
$ws.Range("D150").Value = 21155
$ws.Range("E150").Value = 278330618

$ws.Range("C152").Value = 126010
$ws.Range("E152").Value = 715356771

$ws.Range("C156").Value = 47575
$ws.Range("E156").Value = 142313810

$ws.Range("C168").Value = 284690
$ws.Range("E168").Value = 1201165792

$ws.Range("C169").Value = 562447
$ws.Range("E169").Value = 1283415930

$ws.Range("C170").Value = 366853
$ws.Range("E170").Value = 2838319073

$ws.Range("C171").Value = 114998
$ws.Range("E171").Value = 441773336

$ws.Range("C174").Value = 356935
$ws.Range("E174").Value = 1012747379

$ws.Range("C175").Value = 125349
$ws.Range("E175").Value = 803974249

$ws.Range("C177").Value = 96724
$ws.Range("E177").Value = 174195240

$ws.Range("C179").Value = 235430
$ws.Range("E179").Value = 807980332

$ws.Range("C180").Value = 141375
$ws.Range("E180").Value = 338890651

$ws.Range("C182").Value = 6443
$ws.Range("E182").Value = 12759896

$ws.Range("C199").Value = 4082
$ws.Range("E199").Value = 8780876

$ws.Range("C203").Value = 12769
$ws.Range("E203").Value = 32112207

$ws.Range("C204").Value = 4640
$ws.Range("E204").Value = 11049190

$ws.Range("C205").Value = 10727
$ws.Range("E205").Value = 41578823

$ws.Range("C208").Value = 1515
$ws.Range("E208").Value = 3187444

$ws.Range("C213").Value = 3505
$ws.Range("E213").Value = 10616306

$ws.Range("C214").Value = 6101
$ws.Range("E214").Value = 10896175

$ws.Range("C247").Value = 29417
$ws.Range("E247").Value = 99420663

$ws.Range("C276").Value = 216590
$ws.Range("E276").Value = 1209581281

$ws.Range("C295").Value = 91313
$ws.Range("E295").Value = 552788457

$ws.Range("C298").Value = 11912
$ws.Range("E298").Value = 24040364

$ws.Range("C311").Value = 190830
$ws.Range("E311").Value = 585868780

$ws.Range("C313").Value = 220580
$ws.Range("E313").Value = 1369867687

$ws.Range("C320").Value = 67234
$ws.Range("E320").Value = 124542192

$ws.Range("C322").Value = 81131
$ws.Range("E322").Value = 254162109

$ws.Range("C323").Value = 94713
$ws.Range("E323").Value = 178744833
